$wb = $excel.ActiveWorkbook

# "Generate Report for Handoff": a new handoff was generated for the
# "4d93103c-a919-4d01-b99a-dbe0d34ebbde" file, so its "Latest Handoff
# Datetime" (column D, row 4) is updated to a new timestamp in both the
# zh-cn and de-de status sheets.

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-03-02 09:39:16"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-03-02 09:39:27"
